$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap India and Brasil rows (India now ranks above Brasil) ---
# Row 5 becomes India (with updated stats), Row 6 becomes Brasil (previous Row 5 stats)
$ws.Range("A5").Value = "India"
$ws.Range("A6").Value = "Brasil"

# --- Updated COVID-19 country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B5").Value = 4131690
$ws.Range("C5").Value = 20851
$ws.Range("D5").Value = 3195459
$ws.Range("E5").Value = 865429
$ws.Range("G5").Value = 123
$ws.Range("H5").Value = 70802
$ws.Range("B6").Value = 4123000
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 3296702
$ws.Range("E6").Value = 700068
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 126230
$ws.Range("B18").Value = 320688
$ws.Range("C18").Value = 756
$ws.Range("D18").Value = 296737
$ws.Range("E18").Value = 19870
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 4081
$ws.Range("B24").Value = 251130
$ws.Range("C24").Value = 74
$ws.Range("E24").Value = 15521
$ws.Range("D27").Value = 62227
$ws.Range("E27").Value = 70821
$ws.Range("G27").Value = 35
$ws.Range("H27").Value = 2846
$ws.Range("B51").Value = 60258
$ws.Range("C51").Value = 315
$ws.Range("D51").Value = 42953
$ws.Range("E51").Value = 15465
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 1840
$ws.Range("E55").Value = 3927
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 199
$ws.Range("B67").Value = 37329
$ws.Range("C67").Value = 137
$ws.Range("D67").Value = 34705
$ws.Range("E67").Value = 2076
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 548
$ws.Range("B68").Value = 35103
$ws.Range("C68").Value = 83
$ws.Range("D68").Value = 21230
$ws.Range("E68").Value = 13276
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 597
$ws.Range("B69").Value = 31905
$ws.Range("C69").Value = 56
$ws.Range("D69").Value = 30637
$ws.Range("E69").Value = 544
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 724
$ws.Range("B77").Value = 21560
$ws.Range("C77").Value = 121
$ws.Range("D77").Value = 14709
$ws.Range("E77").Value = 6196
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 655
$ws.Range("B86").Value = 15090
$ws.Range("C86").Value = 92
$ws.Range("D86").Value = 12235
$ws.Range("E86").Value = 2238
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 617
$ws.Range("B89").Value = 12776
$ws.Range("C89").Value = 67
$ws.Range("D89").Value = 11674
$ws.Range("E89").Value = 807
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 295

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 15:53"
